# Insert a new record row at row 583 (a new daily price observation for
# "Betarraga" at "Macroferia Regional de Talca") which pushes all
# subsequent rows (old 583..696) down by one position (new 584..697).
#
# This matches the supplied diff: the dimension grows from A1:R696 to
# A1:R697, and every row from 583 through 696 effectively receives the
# values that used to belong to the row above it, with a brand-new row
# of data occupying the (now vacated) row 583.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 583, shifting rows 583:696 to 584:697
$ws.Rows.Item(583).Insert()

# Populate the newly inserted row 583 with the new record's data
$ws.Range("A583").Value = 5
$ws.Range("B583").Value = "Macroferia Regional de Talca"
$ws.Range("C583").Value = "Maule"
$ws.Range("D583").Value = 45258
$ws.Range("E583").Value = 7
$ws.Range("F583").Value = 100114014
$ws.Range("G583").Value = "Betarraga"
$ws.Range("H583").Value = "Sin especificar"
$ws.Range("I583").Value = "Primera"
$ws.Range("J583").Value = 5000
$ws.Range("K583").Value = 700
$ws.Range("L583").Value = 700
$ws.Range("M583").Value = 700
$ws.Range("N583").Value = "`$/paquete 5 unidades"
$ws.Range("O583").Value = "Región del Maule"
$ws.Range("P583").Value = 140
$ws.Range("Q583").Value = 5
$ws.Range("R583").Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest
# of column D (style index 2 in the original workbook), matching the
# style that was already carried onto the new row by the Insert above.
$ws.Range("D583").NumberFormat = $ws.Range("D584").NumberFormat
